$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 10.100659
$ws.Range("H2").Value = 30.301977
$ws.Range("I2").Value = 0.3328245842863797
$ws.Range("J2").Value = 0.3328245842863797
$ws.Range("M2").Value = 43.73434833333334
$ws.Range("N2").Value = 131.203045
$ws.Range("O2").Value = 0.1998633389969613
$ws.Range("P2").Value = 0.1998633389969613
$ws.Range("Q2").Value = 441.7457391022184
$ws.Range("R2").Value = 3975.711651919965
$ws.Range("S2").Value = 0.06651943271575143
$ws.Range("T2").Value = 0.06651943271575141

# Row 3
$ws.Range("G3").Value = 10.100659
$ws.Range("H3").Value = 30.301977
$ws.Range("I3").Value = 0.3328245842863797
$ws.Range("J3").Value = 0.3328245842863797
$ws.Range("O3").Value = 0.3183113588032023
$ws.Range("P3").Value = 0.3183113588032022
$ws.Range("Q3").Value = 703.5441675538597
$ws.Range("R3").Value = 6331.897507984737
$ws.Range("S3").Value = 0.1059418456673084
$ws.Range("T3").Value = 0.1059418456673084

# Row 4
$ws.Range("G4").Value = 10.100659
$ws.Range("H4").Value = 30.301977
$ws.Range("I4").Value = 0.3328245842863797
$ws.Range("J4").Value = 0.3328245842863797
$ws.Range("M4").Value = 37.39234266666667
$ws.Range("N4").Value = 112.177028
$ws.Range("O4").Value = 0.1708807549004341
$ws.Range("P4").Value = 0.170880754900434
$ws.Range("Q4").Value = 377.6873024871507
$ws.Range("R4").Value = 3399.185722384356
$ws.Range("S4").Value = 0.0568733162122797
$ws.Range("T4").Value = 0.05687331621227969

# Row 5
$ws.Range("G5").Value = 10.100659
$ws.Range("H5").Value = 30.301977
$ws.Range("I5").Value = 0.3328245842863797
$ws.Range("J5").Value = 0.3328245842863797
$ws.Range("M5").Value = 68.04127866666666
$ws.Range("N5").Value = 204.123836
$ws.Range("O5").Value = 0.3109445472994024
$ws.Range("P5").Value = 0.3109445472994024
$ws.Range("Q5").Value = 687.2617537359746
$ws.Range("R5").Value = 6185.355783623771
$ws.Range("S5").Value = 0.1034899896910401
$ws.Range("T5").Value = 0.1034899896910401

# Row 6
$ws.Range("I6").Value = 0.4180918757349671
$ws.Range("J6").Value = 0.4180918757349671
$ws.Range("M6").Value = 43.73434833333334
$ws.Range("N6").Value = 131.203045
$ws.Range("O6").Value = 0.1998633389969613
$ws.Range("P6").Value = 0.1998633389969613
$ws.Range("Q6").Value = 554.9178557683067
$ws.Range("R6").Value = 4994.26070191476
$ws.Range("S6").Value = 0.08356123829189316
$ws.Range("T6").Value = 0.08356123829189313

# Row 7
$ws.Range("I7").Value = 0.4180918757349671
$ws.Range("J7").Value = 0.4180918757349671
$ws.Range("O7").Value = 0.3183113588032023
$ws.Range("P7").Value = 0.3183113588032022
$ws.Range("S7").Value = 0.133083393069777
$ws.Range("T7").Value = 0.133083393069777

# Row 8
$ws.Range("I8").Value = 0.4180918757349671
$ws.Range("J8").Value = 0.4180918757349671
$ws.Range("M8").Value = 37.39234266666667
$ws.Range("N8").Value = 112.177028
$ws.Range("O8").Value = 0.1708807549004341
$ws.Range("P8").Value = 0.170880754900434
$ws.Range("Q8").Value = 474.4481032755094
$ws.Range("R8").Value = 4270.032929479585
$ws.Range("S8").Value = 0.07144385534332964
$ws.Range("T8").Value = 0.07144385534332963

# Row 9
$ws.Range("I9").Value = 0.4180918757349671
$ws.Range("J9").Value = 0.4180918757349671
$ws.Range("M9").Value = 68.04127866666666
$ws.Range("N9").Value = 204.123836
$ws.Range("O9").Value = 0.3109445472994024
$ws.Range("P9").Value = 0.3109445472994024
$ws.Range("Q9").Value = 863.3333272434452
$ws.Range("R9").Value = 7769.999945191008
$ws.Range("S9").Value = 0.1300033890299674
$ws.Range("T9").Value = 0.1300033890299674

# Row 10
$ws.Range("G10").Value = 4.721016333333334
$ws.Range("H10").Value = 14.163049
$ws.Range("I10").Value = 0.1555611667071302
$ws.Range("J10").Value = 0.1555611667071302
$ws.Range("M10").Value = 43.73434833333334
$ws.Range("N10").Value = 131.203045
$ws.Range("O10").Value = 0.1998633389969613
$ws.Range("P10").Value = 0.1998633389969613
$ws.Range("Q10").Value = 206.4705728093562
$ws.Range("R10").Value = 1858.235155284205
$ws.Range("S10").Value = 0.03109097419634999
$ws.Range("T10").Value = 0.03109097419634997

# Row 11
$ws.Range("G11").Value = 4.721016333333334
$ws.Range("H11").Value = 14.163049
$ws.Range("I11").Value = 0.1555611667071302
$ws.Range("J11").Value = 0.1555611667071302
$ws.Range("O11").Value = 0.3183113588032023
$ws.Range("P11").Value = 0.3183113588032022
$ws.Range("Q11").Value = 328.8343370707966
$ws.Range("R11").Value = 2959.509033637169
$ws.Range("S11").Value = 0.0495168863515581
$ws.Range("T11").Value = 0.04951688635155808

# Row 12
$ws.Range("G12").Value = 4.721016333333334
$ws.Range("H12").Value = 14.163049
$ws.Range("I12").Value = 0.1555611667071302
$ws.Range("J12").Value = 0.1555611667071302
$ws.Range("M12").Value = 37.39234266666667
$ws.Range("N12").Value = 112.177028
$ws.Range("O12").Value = 0.1708807549004341
$ws.Range("P12").Value = 0.170880754900434
$ws.Range("Q12").Value = 176.5298604709303
$ws.Range("R12").Value = 1588.768744238372
$ws.Range("S12").Value = 0.02658240960010668
$ws.Range("T12").Value = 0.02658240960010668

# Row 13
$ws.Range("G13").Value = 4.721016333333334
$ws.Range("H13").Value = 14.163049
$ws.Range("I13").Value = 0.1555611667071302
$ws.Range("J13").Value = 0.1555611667071302
$ws.Range("M13").Value = 68.04127866666666
$ws.Range("N13").Value = 204.123836
$ws.Range("O13").Value = 0.3109445472994024
$ws.Range("P13").Value = 0.3109445472994024
$ws.Range("Q13").Value = 321.2239879262182
$ws.Range("R13").Value = 2891.015891335964
$ws.Range("S13").Value = 0.04837089655911549
$ws.Range("T13").Value = 0.04837089655911547

# Row 14
$ws.Range("G14").Value = 2.838244666666667
$ws.Range("H14").Value = 8.514734000000001
$ws.Range("I14").Value = 0.09352237327152295
$ws.Range("J14").Value = 0.09352237327152294
$ws.Range("M14").Value = 43.73434833333334
$ws.Range("N14").Value = 131.203045
$ws.Range("O14").Value = 0.1998633389969613
$ws.Range("P14").Value = 0.1998633389969613
$ws.Range("Q14").Value = 124.1287809072256
$ws.Range("R14").Value = 1117.15902816503
$ws.Range("S14").Value = 0.01869169379296675
$ws.Range("T14").Value = 0.01869169379296674

# Row 15
$ws.Range("G15").Value = 2.838244666666667
$ws.Range("H15").Value = 8.514734000000001
$ws.Range("I15").Value = 0.09352237327152295
$ws.Range("J15").Value = 0.09352237327152294
$ws.Range("O15").Value = 0.3183113588032023
$ws.Range("P15").Value = 0.3183113588032022
$ws.Range("Q15").Value = 197.6930892651838
$ws.Range("R15").Value = 1779.237803386654
$ws.Range("S15").Value = 0.02976923371455876
$ws.Range("T15").Value = 0.02976923371455874

# Row 16
$ws.Range("G16").Value = 2.838244666666667
$ws.Range("H16").Value = 8.514734000000001
$ws.Range("I16").Value = 0.09352237327152295
$ws.Range("J16").Value = 0.09352237327152294
$ws.Range("M16").Value = 37.39234266666667
$ws.Range("N16").Value = 112.177028
$ws.Range("O16").Value = 0.1708807549004341
$ws.Range("P16").Value = 0.170880754900434
$ws.Range("Q16").Value = 106.1286171478391
$ws.Range("R16").Value = 955.1575543305521
$ws.Range("S16").Value = 0.01598117374471802
$ws.Range("T16").Value = 0.01598117374471801

# Row 17
$ws.Range("G17").Value = 2.838244666666667
$ws.Range("H17").Value = 8.514734000000001
$ws.Range("I17").Value = 0.09352237327152295
$ws.Range("J17").Value = 0.09352237327152294
$ws.Range("M17").Value = 68.04127866666666
$ws.Range("N17").Value = 204.123836
$ws.Range("O17").Value = 0.3109445472994024
$ws.Range("P17").Value = 0.3109445472994024
$ws.Range("Q17").Value = 193.1177962888471
$ws.Range("R17").Value = 1738.060166599624
$ws.Range("S17").Value = 0.02908027201927944
$ws.Range("T17").Value = 0.02908027201927943

Write-Host "Updated Hbegf-Cd9 TPM values"